$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 509, shifting the existing rows 509:541 down to 510:542.
$ws.Rows.Item(509).Insert()

# Populate the newly inserted row 509. It is a new weekly price entry for
# "Macroferia Regional de Talca" / Acelga, duplicating the previous latest
# entry's price data but dated one period later (2023-12-05, serial 45265).
$ws.Cells.Item(509, 1).Value = 5
$ws.Cells.Item(509, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(509, 3).Value = "Maule"
$ws.Cells.Item(509, 4).Value = 45265
$ws.Cells.Item(509, 5).Value = 7
$ws.Cells.Item(509, 6).Value = 100112009
$ws.Cells.Item(509, 7).Value = "Acelga"
$ws.Cells.Item(509, 8).Value = "Sin especificar"
$ws.Cells.Item(509, 9).Value = "Primera"
$ws.Cells.Item(509, 10).Value = 500
$ws.Cells.Item(509, 11).Value = 2500
$ws.Cells.Item(509, 12).Value = 2500
$ws.Cells.Item(509, 13).Value = 2500
$ws.Cells.Item(509, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(509, 15).Value = "Región del Maule"
$ws.Cells.Item(509, 16).Value = 625
$ws.Cells.Item(509, 17).Value = 4
$ws.Cells.Item(509, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(509, 4).NumberFormat = $ws.Cells.Item(510, 4).NumberFormat
